# MVP of irt0 model class
# Adds an 11th item column ("J") to the item-response matrix:
#  - new shared string "J" for the row-11 label
#  - new column L holding the per-row item count (11 in the header row,
#    1 for every data row)
#  - new row 11 for item "J" with all-zero responses except the new
#    column L marker
#  - refreshes the active selection the way the original author left it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new column L (header value 11, then 1 for every data row 2-11) ---
$ws.Range("L1").Value = 11
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("L9").Value = 1
$ws.Range("L10").Value = 1
$ws.Range("L11").Value = 1

# --- new row 11 for item "J" ---
$ws.Range("A11").Value = "J"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0

# --- restore the author's selection ---
$ws.Range("J15").Select() | Out-Null
